# Updates the cryptos price/volume table with refreshed figures.
# D-column values are prefixed with a leading apostrophe so Excel keeps
# them as literal text (matching the workbook's original inlineStr cells)
# instead of re-parsing them as numbers and mangling formatting such as
# trailing zeros (e.g. "323.60" -> 323.6) or thousand-separated price
# strings (e.g. "27.366.71").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.366.71'
$ws.Range("E2").Value = '  -3.42%  '

$ws.Range("D3").Value = '''1.859.28'
$ws.Range("E3").Value = '  -4.23%  '

$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -1.21%  '

$ws.Range("D5").Value = '''323.60'
$ws.Range("E5").Value = '  +0.78%  '

$ws.Range("D6").Value = '''1.001'

$ws.Range("D7").Value = '''0.4537'
$ws.Range("E7").Value = '  -4.89%  '

$ws.Range("D8").Value = '''0.3866'
$ws.Range("E8").Value = '  -5.03%  '

$ws.Range("D9").Value = '''48.43'
$ws.Range("E9").Value = '  -9.43%  '

$ws.Range("D10").Value = '''0.07916'

$ws.Range("E11").Value = '  -3.40%  '

$ws.Range("E12").Value = '  -4.09%  '

$ws.Range("D13").Value = '''1.863.44'
$ws.Range("E13").Value = '  -7.05%  '

$ws.Range("D14").Value = '''5.905'
$ws.Range("E14").Value = '  -3.69%  '

$ws.Range("D15").Value = '''7.125'
$ws.Range("E15").Value = '  -5.60%  '

$ws.Range("D16").Value = '''1.001'
$ws.Range("E16").Value = '  -1.30%  '

$ws.Range("E17").Value = '  -3.77%  '

$ws.Range("D18").Value = '''85.84'
$ws.Range("E18").Value = '  -4.70%  '

$ws.Range("D19").Value = '''0.06512'
$ws.Range("E19").Value = '  -1.99%  '

$ws.Range("D20").Value = '''17.06'
$ws.Range("E20").Value = '  -6.80%  '

$ws.Range("E21").Value = '  -1.06%  '

$ws.Range("D22").Value = '''5.530'
$ws.Range("E22").Value = '  -5.04%  '

$ws.Range("D23").Value = '''27.368.14'
$ws.Range("E23").Value = '  -3.63%  '

$ws.Range("D24").Value = '''10.85'
$ws.Range("E24").Value = '  -5.07%  '

$ws.Range("D25").Value = '''2.278'
$ws.Range("E25").Value = '  -0.89%  '

$ws.Range("D26").Value = '''2.086.23'
$ws.Range("E26").Value = '  -6.52%  '

$ws.Range("D27").Value = '''153.69'
$ws.Range("E27").Value = '  -1.33%  '

$ws.Range("D28").Value = '''19.75'
$ws.Range("E28").Value = '  -2.45%  '

$ws.Range("D29").Value = '''2.073'
$ws.Range("E29").Value = '  -4.61%  '

$ws.Range("D30").Value = '''5.434'
$ws.Range("E30").Value = '  -6.24%  '

$ws.Range("D31").Value = '''120.71'
$ws.Range("E31").Value = '  -2.58%  '

$ws.Range("D32").Value = '''1.485'
$ws.Range("E32").Value = '  +2.77%  '

$ws.Range("D33").Value = '''0.09298'
$ws.Range("E33").Value = '  -3.43%  '

$ws.Range("D34").Value = '''0.9360'
$ws.Range("E34").Value = '  -4.63%  '

$ws.Range("E35").Value = '  -2.64%  '

$ws.Range("D36").Value = '''5.257'
$ws.Range("E36").Value = '  -6.32%  '

$ws.Range("D37").Value = '''0.02238'
$ws.Range("E37").Value = '  -3.79%  '

$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '''0.05989'
$ws.Range("E38").Value = '  -3.01%  '

$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = '''1.220'
$ws.Range("E39").Value = '  -1.97%  '

$ws.Range("D40").Value = '''8.223'
$ws.Range("E40").Value = '  -10.29%  '

$ws.Range("E41").Value = '  -1.06%  '

$ws.Range("D42").Value = '''0.5909'
$ws.Range("E42").Value = '  -4.70%  '

$ws.Range("D43").Value = '''0.1888'
$ws.Range("E43").Value = '  -1.42%  '

$ws.Range("D44").Value = '''10.12'
$ws.Range("E44").Value = '  -9.51%  '

$ws.Range("D45").Value = '''1.276'
$ws.Range("E45").Value = '  -3.92%  '

$ws.Range("D46").Value = '''0.5619'
$ws.Range("E46").Value = '  -5.37%  '

$ws.Range("D47").Value = '''12.04'
$ws.Range("E47").Value = '  -6.02%  '

$ws.Range("D48").Value = '''3.365'
$ws.Range("E48").Value = '  -1.12%  '

$ws.Range("D49").Value = '''1.922'
$ws.Range("E49").Value = '  -6.48%  '

$ws.Range("D50").Value = '''0.06772'
$ws.Range("E50").Value = '  -0.38%  '

$ws.Range("D51").Value = '''108.17'
$ws.Range("E51").Value = '  -1.93%  '
